$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 62
$ws.Range("H62").Value = 1054
$ws.Range("I62").Value = 998.25
$ws.Range("J62").Value = 1500
$ws.Range("K62").Value = 998.25
$ws.Range("L62").Value = 1500
$ws.Range("M62").Value = -374.25
$ws.Range("N62").Value = -2748
# Row 65
$ws.Range("H65").Value = 1054
$ws.Range("I65").Value = 998.25
$ws.Range("J65").Value = 1500
$ws.Range("K65").Value = 4991.25
$ws.Range("L65").Value = 7500
$ws.Range("M65").Value = -1871.25
$ws.Range("N65").Value = -13740
# Row 127
$ws.Range("H127").Value = 44852.695
$ws.Range("I127").Value = 297.25
$ws.Range("J127").Value = 68615.60000000001
$ws.Range("K127").Value = 891.75
$ws.Range("L127").Value = 205846.8
$ws.Range("M127").Value = 4068.25
$ws.Range("N127").Value = -215766.8
# Row 135
$ws.Range("H135").Value = 351.47827
$ws.Range("I135").Value = 320.36365
$ws.Range("J135").Value = 1036
$ws.Range("K135").Value = 2883.27285
$ws.Range("L135").Value = 9324
$ws.Range("M135").Value = -348.2728500000003
$ws.Range("N135").Value = -14394

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 9870.437
$ws.Range("I32").Value = 5594.2856
$ws.Range("J32").Value = 23685.691
$ws.Range("K32").Value = 5594.2856
$ws.Range("L32").Value = 23685.691
$ws.Range("M32").Value = -5307.2856
$ws.Range("N32").Value = -24259.691
# Row 45
$ws.Range("H45").Value = 1362.6666
$ws.Range("I45").Value = 1206.5
$ws.Range("J45").Value = 1675
$ws.Range("K45").Value = 1206.5
$ws.Range("L45").Value = 1675
$ws.Range("M45").Value = -829.5
$ws.Range("N45").Value = -2429
# Row 74
$ws.Range("H74").Value = 1591.1765
$ws.Range("I74").Value = 1393.5714
$ws.Range("J74").Value = 2513.3333
$ws.Range("K74").Value = 1393.5714
$ws.Range("L74").Value = 2513.3333
$ws.Range("M74").Value = -519.5714
$ws.Range("N74").Value = -4261.3333
# Row 77
$ws.Range("H77").Value = 1591.1765
$ws.Range("I77").Value = 1393.5714
$ws.Range("J77").Value = 2513.3333
$ws.Range("K77").Value = 6967.857
$ws.Range("L77").Value = 12566.6665
$ws.Range("M77").Value = -2599.857
$ws.Range("N77").Value = -21302.6665
# Row 110
$ws.Range("H110").Value = 4799.8
$ws.Range("I110").Value = 999.5
$ws.Range("J110").Value = 7333.3335
$ws.Range("K110").Value = 999.5
$ws.Range("L110").Value = 7333.3335
$ws.Range("M110").Value = 1045.5
$ws.Range("N110").Value = -11423.3335
# Row 132
$ws.Range("H132").Value = 2275.111
$ws.Range("I132").Value = 1611.9524
$ws.Range("J132").Value = 3203.5334
$ws.Range("K132").Value = 4835.857199999999
$ws.Range("L132").Value = 9610.600199999999
$ws.Range("M132").Value = -2305.857199999999
$ws.Range("N132").Value = -14670.6002

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 26
$ws.Range("H26").Value = 45117.75
$ws.Range("I26").Value = 30235.5
$ws.Range("J26").Value = 60000
$ws.Range("K26").Value = 30235.5
$ws.Range("L26").Value = 60000
$ws.Range("M26").Value = -29943.5
$ws.Range("N26").Value = -60584
# Row 40
$ws.Range("H40").Value = 44350
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 44350
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 44350
$ws.Range("N40").Value = -44880
# Row 93
$ws.Range("H93").Value = 47500
$ws.Range("I93").Value = 0
$ws.Range("J93").Value = 47500
$ws.Range("K93").Value = 0
$ws.Range("L93").Value = 47500
$ws.Range("N93").Value = -51244
# Row 95
$ws.Range("H95").Value = 19680.445
$ws.Range("I95").Value = 0
$ws.Range("J95").Value = 19680.445
$ws.Range("K95").Value = 0
$ws.Range("L95").Value = 19680.445
$ws.Range("N95").Value = -25172.445
# Row 96
$ws.Range("H96").Value = 40025.6
$ws.Range("I96").Value = 10714
$ws.Range("J96").Value = 59566.668
$ws.Range("K96").Value = 10714
$ws.Range("L96").Value = 59566.668
$ws.Range("M96").Value = -7968
$ws.Range("N96").Value = -65058.668

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 37
$ws.Range("H37").Value = 42500
$ws.Range("I37").Value = 40000
$ws.Range("J37").Value = 45000
$ws.Range("K37").Value = 40000
$ws.Range("L37").Value = 45000
$ws.Range("M37").Value = -39893
$ws.Range("N37").Value = -45214
# Row 99
$ws.Range("H99").Value = 9270.929
$ws.Range("I99").Value = 11579.3
$ws.Range("J99").Value = 3500
$ws.Range("K99").Value = 11579.3
$ws.Range("L99").Value = 3500
$ws.Range("M99").Value = -10081.3
$ws.Range("N99").Value = -6496
# Row 126
$ws.Range("H126").Value = 9270.929
$ws.Range("I126").Value = 11579.3
$ws.Range("J126").Value = 3500
$ws.Range("K126").Value = 34737.89999999999
$ws.Range("L126").Value = 10500
$ws.Range("M126").Value = -32267.89999999999
$ws.Range("N126").Value = -15440
# Row 132
$ws.Range("H132").Value = 1460.9286
$ws.Range("I132").Value = 1011.25806
$ws.Range("J132").Value = 2728.182
$ws.Range("K132").Value = 3033.77418
$ws.Range("L132").Value = 8184.545999999999
$ws.Range("M132").Value = -503.7741799999999
$ws.Range("N132").Value = -13244.546

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 69
$ws.Range("H69").Value = 1040
$ws.Range("I69").Value = 800
$ws.Range("J69").Value = 1100
$ws.Range("K69").Value = 2400
$ws.Range("L69").Value = 3300
$ws.Range("M69").Value = -1589
$ws.Range("N69").Value = -4922
# Row 72
$ws.Range("H72").Value = 1040
$ws.Range("I72").Value = 800
$ws.Range("J72").Value = 1100
$ws.Range("K72").Value = 7200
$ws.Range("L72").Value = 9900
$ws.Range("M72").Value = -3144
$ws.Range("N72").Value = -18012
# Row 113
$ws.Range("H113").Value = 7936947
$ws.Range("I113").Value = 438.33334
$ws.Range("J113").Value = 23809964
$ws.Range("K113").Value = 1315.00002
$ws.Range("L113").Value = 71429892
$ws.Range("M113").Value = 854.9999800000001
$ws.Range("N113").Value = -71434232
# Row 121
$ws.Range("H121").Value = 3795.8823
$ws.Range("I121").Value = 315
$ws.Range("J121").Value = 4013.4375
$ws.Range("K121").Value = 945
$ws.Range("L121").Value = 12040.3125
$ws.Range("M121").Value = 365
$ws.Range("N121").Value = -14660.3125
# Row 136
$ws.Range("H136").Value = 5832.727
$ws.Range("I136").Value = 2610
$ws.Range("J136").Value = 9700
$ws.Range("K136").Value = 7830
$ws.Range("L136").Value = 29100
$ws.Range("M136").Value = -2730
$ws.Range("N136").Value = -39300

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 41
$ws.Range("H41").Value = 14485.714
$ws.Range("I41").Value = 4500
$ws.Range("J41").Value = 21975
$ws.Range("K41").Value = 4500
$ws.Range("L41").Value = 21975
$ws.Range("M41").Value = -4145
$ws.Range("N41").Value = -22685

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7
$ws.Range("H7").Value = 2601
$ws.Range("I7").Value = 2000
$ws.Range("J7").Value = 3502.5
$ws.Range("K7").Value = 2000
$ws.Range("L7").Value = 3502.5
$ws.Range("M7").Value = -1888
$ws.Range("N7").Value = -3726.5
# Row 22
$ws.Range("H22").Value = 566.3
$ws.Range("I22").Value = 421.33334
$ws.Range("J22").Value = 628.4286
$ws.Range("K22").Value = 421.33334
$ws.Range("L22").Value = 628.4286
$ws.Range("M22").Value = -126.33334
$ws.Range("N22").Value = -1218.4286
# Row 27
$ws.Range("H27").Value = 566.3
$ws.Range("I27").Value = 421.33334
$ws.Range("J27").Value = 628.4286
$ws.Range("K27").Value = 421.33334
$ws.Range("L27").Value = 628.4286
$ws.Range("M27").Value = -314.33334
$ws.Range("N27").Value = -842.4286
# Row 94
$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 0
$ws.Range("N94").Value = ""
# Row 122
$ws.Range("H122").Value = 6469
$ws.Range("I122").Value = 7417
$ws.Range("J122").Value = 3625
$ws.Range("K122").Value = 22251
$ws.Range("L122").Value = 10875
$ws.Range("M122").Value = -19801
$ws.Range("N122").Value = -15775
# Row 126
$ws.Range("H126").Value = 2601
$ws.Range("I126").Value = 2000
$ws.Range("J126").Value = 3502.5
$ws.Range("K126").Value = 6000
$ws.Range("L126").Value = 10507.5
$ws.Range("M126").Value = -3530
$ws.Range("N126").Value = -15447.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 94
$ws.Range("H94").Value = 33500
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 33500
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 33500
$ws.Range("N94").Value = -35302
# Row 122
$ws.Range("H122").Value = 1337.375
$ws.Range("I122").Value = 1159.8
$ws.Range("J122").Value = 1633.3334
$ws.Range("K122").Value = 3479.4
$ws.Range("L122").Value = 4900.0002
$ws.Range("M122").Value = -1029.4
$ws.Range("N122").Value = -9800.0002
# Row 136
$ws.Range("H136").Value = 2248
$ws.Range("I136").Value = 2248
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 6744
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -4194
$ws.Range("N136").Value = ""
